$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 3576
$ws.Range("I52").Value = 229
$ws.Range("J52").Value = 5249.5
$ws.Range("K52").Value = 687
$ws.Range("L52").Value = 15748.5
$ws.Range("M52").Value = -527
$ws.Range("N52").Value = -16068.5
$ws.Range("H101").Value = 14291024
$ws.Range("I101").Value = 23814852
$ws.Range("K101").Value = 71444556
$ws.Range("M101").Value = -71442934
$ws.Range("H141").Value = 4927.9033
$ws.Range("I141").Value = 4732.2593
$ws.Range("J141").Value = 6248.5
$ws.Range("K141").Value = 14196.7779
$ws.Range("L141").Value = 18745.5
$ws.Range("M141").Value = -9016.777899999999
$ws.Range("N141").Value = -29105.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("H32").Value = 7229.127
$ws.Range("I32").Value = 7090.5835
$ws.Range("K32").Value = 7090.5835
$ws.Range("M32").Value = -6803.5835
$ws.Range("H61").Value = 13027
$ws.Range("I61").Value = 18140.182
$ws.Range("J61").Value = 6777.5557
$ws.Range("K61").Value = 18140.182
$ws.Range("L61").Value = 6777.5557
$ws.Range("M61").Value = -17928.182
$ws.Range("N61").Value = -7201.5557
$ws.Range("H136").Value = 13027
$ws.Range("I136").Value = 18140.182
$ws.Range("J136").Value = 6777.5557
$ws.Range("K136").Value = 54420.546
$ws.Range("L136").Value = 20332.6671
$ws.Range("M136").Value = -51870.546
$ws.Range("N136").Value = -25432.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 2256.0557
$ws.Range("I107").Value = 2271.1177
$ws.Range("K107").Value = 2271.1177
$ws.Range("M107").Value = -351.1176999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1312.375
$ws.Range("I16").Value = 1237.25
$ws.Range("K16").Value = 1237.25
$ws.Range("M16").Value = -950.25
$ws.Range("H31").Value = 6666.2256
$ws.Range("I31").Value = 7393.909
$ws.Range("J31").Value = 4887.4443
$ws.Range("K31").Value = 7393.909
$ws.Range("L31").Value = 4887.4443
$ws.Range("M31").Value = -7098.909
$ws.Range("N31").Value = -5477.4443
$ws.Range("H34").Value = 6666.2256
$ws.Range("I34").Value = 7393.909
$ws.Range("J34").Value = 4887.4443
$ws.Range("K34").Value = 7393.909
$ws.Range("L34").Value = 4887.4443
$ws.Range("M34").Value = -7191.909
$ws.Range("N34").Value = -5291.4443
$ws.Range("H113").Value = 1312.375
$ws.Range("I113").Value = 1237.25
$ws.Range("K113").Value = 1237.25
$ws.Range("M113").Value = 932.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 69516.53999999999
$ws.Range("I80").Value = 8350.666999999999
$ws.Range("J80").Value = 87866.3
$ws.Range("K80").Value = 25052.001
$ws.Range("L80").Value = 263598.9
$ws.Range("M80").Value = -24116.001
$ws.Range("N80").Value = -265470.9
$ws.Range("H82").Value = 6722
$ws.Range("I82").Value = 972.5
$ws.Range("K82").Value = 2917.5
$ws.Range("M82").Value = -2511.5
$ws.Range("H83").Value = 69516.53999999999
$ws.Range("I83").Value = 8350.666999999999
$ws.Range("J83").Value = 87866.3
$ws.Range("K83").Value = 75156.003
$ws.Range("L83").Value = 790796.7000000001
$ws.Range("M83").Value = -70476.003
$ws.Range("N83").Value = -800156.7000000001
$ws.Range("H85").Value = 6722
$ws.Range("I85").Value = 972.5
$ws.Range("K85").Value = 2917.5
$ws.Range("M85").Value = -1513.5
$ws.Range("H87").Value = 19570.908
$ws.Range("I87").Value = 18600
$ws.Range("J87").Value = 19786.666
$ws.Range("K87").Value = 55800
$ws.Range("L87").Value = 59359.99800000001
$ws.Range("M87").Value = -54552
$ws.Range("N87").Value = -61855.99800000001
$ws.Range("H90").Value = 19570.908
$ws.Range("I90").Value = 18600
$ws.Range("J90").Value = 19786.666
$ws.Range("K90").Value = 167400
$ws.Range("L90").Value = 178079.994
$ws.Range("M90").Value = -161160
$ws.Range("N90").Value = -190559.994

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 5061.3125
$ws.Range("I97").Value = 5177.3213
$ws.Range("K97").Value = 5177.3213
$ws.Range("M97").Value = -4681.3213
$ws.Range("H132").Value = 3457
$ws.Range("I132").Value = 3555.2974
$ws.Range("J132").Value = 2547.75
$ws.Range("K132").Value = 10665.8922
$ws.Range("L132").Value = 7643.25
$ws.Range("M132").Value = -8135.8922
$ws.Range("N132").Value = -12703.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H64").Value = 59000
$ws.Range("J64").Value = 59000
$ws.Range("L64").Value = 59000
$ws.Range("N64").Value = -59450
$ws.Range("H67").Value = 59000
$ws.Range("J67").Value = 59000
$ws.Range("L67").Value = 59000
$ws.Range("N67").Value = -60560
$ws.Range("H122").Value = 7338.0557
$ws.Range("I122").Value = 6923.75
$ws.Range("K122").Value = 20771.25
$ws.Range("M122").Value = -18321.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 385166.66
$ws.Range("I62").Value = 850723
$ws.Range("J62").Value = 12721.6
$ws.Range("K62").Value = 850723
$ws.Range("L62").Value = 12721.6
$ws.Range("M62").Value = -850099
$ws.Range("N62").Value = -13969.6
$ws.Range("H63").Value = 0
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H65").Value = 385166.66
$ws.Range("I65").Value = 850723
$ws.Range("J65").Value = 12721.6
$ws.Range("K65").Value = 4253615
$ws.Range("L65").Value = 63608
$ws.Range("M65").Value = -4250495
$ws.Range("N65").Value = -69848
$ws.Range("H66").Value = 0
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H75").Value = 39581.668
$ws.Range("J75").Value = 41998
$ws.Range("L75").Value = 41998
$ws.Range("N75").Value = -43870
$ws.Range("H78").Value = 39581.668
$ws.Range("J78").Value = 41998
$ws.Range("L78").Value = 125994
$ws.Range("N78").Value = -135354
$ws.Range("H132").Value = 13099.177
$ws.Range("I132").Value = 18167.21
$ws.Range("J132").Value = 6679.6665
$ws.Range("K132").Value = 54501.63
$ws.Range("L132").Value = 20038.9995
$ws.Range("M132").Value = -51971.63
$ws.Range("N132").Value = -25098.9995
